$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H29").Value = 1499
$ws.Range("I29").Value = 998
$ws.Range("J29").Value = 2000
$ws.Range("K29").Value = 2994
$ws.Range("L29").Value = 6000
$ws.Range("M29").Value = -2713
$ws.Range("N29").Value = -6562
$ws.Range("H38").Value = 2080.5454
$ws.Range("I38").Value = 931.7778
$ws.Range("J38").Value = 7250
$ws.Range("K38").Value = 2795.3334
$ws.Range("L38").Value = 21750
$ws.Range("M38").Value = -2423.3334
$ws.Range("N38").Value = -22494
$ws.Range("H58").Value = 1725.3334
$ws.Range("I58").Value = 969.5714
$ws.Range("J58").Value = 2783.4
$ws.Range("K58").Value = 2908.7142
$ws.Range("L58").Value = 8350.200000000001
$ws.Range("M58").Value = -2758.7142
$ws.Range("N58").Value = -8650.200000000001
$ws.Range("H76").Value = 4189.6
$ws.Range("I76").Value = 4237
$ws.Range("J76").Value = 4000
$ws.Range("K76").Value = 4237
$ws.Range("L76").Value = 4000
$ws.Range("M76").Value = -3922
$ws.Range("N76").Value = -4630
$ws.Range("H79").Value = 4189.6
$ws.Range("I79").Value = 4237
$ws.Range("J79").Value = 4000
$ws.Range("K79").Value = 4237
$ws.Range("L79").Value = 4000
$ws.Range("M79").Value = -3145
$ws.Range("N79").Value = -6184
$ws.Range("H98").Value = 37912
$ws.Range("I98").Value = 61614.668
$ws.Range("J98").Value = 11246.5
$ws.Range("K98").Value = 61614.668
$ws.Range("L98").Value = 11246.5
$ws.Range("M98").Value = -60116.668
$ws.Range("N98").Value = -14242.5
$ws.Range("H116").Value = 790271.5600000001
$ws.Range("I116").Value = 1416310.1
$ws.Range("J116").Value = 7723.375
$ws.Range("K116").Value = 1416310.1
$ws.Range("L116").Value = 7723.375
$ws.Range("M116").Value = -1412868.1
$ws.Range("N116").Value = -14607.375
$ws.Range("H122").Value = 37912
$ws.Range("I122").Value = 61614.668
$ws.Range("J122").Value = 11246.5
$ws.Range("K122").Value = 184844.004
$ws.Range("L122").Value = 33739.5
$ws.Range("M122").Value = -182394.004
$ws.Range("N122").Value = -38639.5
$ws.Range("H132").Value = 2500.951
$ws.Range("I132").Value = 2682.102
$ws.Range("J132").Value = 1761.25
$ws.Range("K132").Value = 8046.306
$ws.Range("L132").Value = 5283.75
$ws.Range("M132").Value = -5516.306
$ws.Range("N132").Value = -10343.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 502250.75
$ws.Range("I11").Value = 668334.3
$ws.Range("J11").Value = 4000
$ws.Range("K11").Value = 668334.3
$ws.Range("L11").Value = 4000
$ws.Range("M11").Value = -668190.3
$ws.Range("N11").Value = -4288
$ws.Range("H32").Value = 3721.8857
$ws.Range("I32").Value = 3769.9678
$ws.Range("J32").Value = 3349.25
$ws.Range("K32").Value = 3769.9678
$ws.Range("L32").Value = 3349.25
$ws.Range("M32").Value = -3482.9678
$ws.Range("N32").Value = -3923.25
$ws.Range("H61").Value = 5957.5
$ws.Range("I61").Value = 5957.5
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 5957.5
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -5745.5
$ws.Range("H132").Value = 3545.6843
$ws.Range("I132").Value = 2894
$ws.Range("J132").Value = 4957.6665
$ws.Range("K132").Value = 8682
$ws.Range("L132").Value = 14872.9995
$ws.Range("M132").Value = -6152
$ws.Range("N132").Value = -19932.9995
$ws.Range("H136").Value = 5957.5
$ws.Range("I136").Value = 5957.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 17872.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -15322.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2265.1333
$ws.Range("I20").Value = 1336.0555
$ws.Range("J20").Value = 3658.75
$ws.Range("K20").Value = 1336.0555
$ws.Range("L20").Value = 3658.75
$ws.Range("M20").Value = -1089.0555
$ws.Range("N20").Value = -4152.75
$ws.Range("H75").Value = 29750
$ws.Range("I75").Value = 14500
$ws.Range("J75").Value = 45000
$ws.Range("K75").Value = 14500
$ws.Range("L75").Value = 45000
$ws.Range("M75").Value = -13564
$ws.Range("N75").Value = -46872
$ws.Range("H78").Value = 29750
$ws.Range("I78").Value = 14500
$ws.Range("J78").Value = 45000
$ws.Range("K78").Value = 43500
$ws.Range("L78").Value = 135000
$ws.Range("M78").Value = -38820
$ws.Range("N78").Value = -144360
$ws.Range("H82").Value = 70796.5
$ws.Range("I82").Value = 27499.5
$ws.Range("J82").Value = 92445
$ws.Range("K82").Value = 27499.5
$ws.Range("L82").Value = 92445
$ws.Range("M82").Value = -27116.5
$ws.Range("N82").Value = -93211
$ws.Range("H85").Value = 70796.5
$ws.Range("I85").Value = 27499.5
$ws.Range("J85").Value = 92445
$ws.Range("K85").Value = 27499.5
$ws.Range("L85").Value = 92445
$ws.Range("M85").Value = -26173.5
$ws.Range("N85").Value = -95097
$ws.Range("H97").Value = 22498.5
$ws.Range("I97").Value = 14997
$ws.Range("J97").Value = 30000
$ws.Range("K97").Value = 14997
$ws.Range("L97").Value = 30000
$ws.Range("M97").Value = -14006
$ws.Range("N97").Value = -31982
$ws.Range("H105").Value = 2423.0557
$ws.Range("I105").Value = 2104
$ws.Range("J105").Value = 3539.75
$ws.Range("K105").Value = 2104
$ws.Range("L105").Value = 3539.75
$ws.Range("M105").Value = -357
$ws.Range("N105").Value = -7033.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 10000
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 10000
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 10000
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -10340
$ws.Range("H31").Value = 6264.4375
$ws.Range("I31").Value = 5438.5
$ws.Range("J31").Value = 7090.375
$ws.Range("K31").Value = 5438.5
$ws.Range("L31").Value = 7090.375
$ws.Range("M31").Value = -5143.5
$ws.Range("N31").Value = -7680.375
$ws.Range("H34").Value = 6264.4375
$ws.Range("I34").Value = 5438.5
$ws.Range("J34").Value = 7090.375
$ws.Range("K34").Value = 5438.5
$ws.Range("L34").Value = 7090.375
$ws.Range("M34").Value = -5236.5
$ws.Range("N34").Value = -7494.375
$ws.Range("H120").Value = 51900
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 51900
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 51900
$ws.Range("N120").Value = -59158

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 10
$ws.Range("I7").Value = 10
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 30
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 82
$ws.Range("H9").Value = 11160.833
$ws.Range("I9").Value = 16457.625
$ws.Range("J9").Value = 567.25
$ws.Range("K9").Value = 49372.875
$ws.Range("L9").Value = 1701.75
$ws.Range("M9").Value = -49148.875
$ws.Range("N9").Value = -2149.75
$ws.Range("H11").Value = 250274.75
$ws.Range("I11").Value = 49.5
$ws.Range("J11").Value = 500500
$ws.Range("K11").Value = 148.5
$ws.Range("L11").Value = 1501500
$ws.Range("M11").Value = -8.5
$ws.Range("N11").Value = -1501780
$ws.Range("H37").Value = 62628.332
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 62628.332
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 187884.996
$ws.Range("N37").Value = -188108.996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8919.15
$ws.Range("I70").Value = 8367.846
$ws.Range("J70").Value = 9943
$ws.Range("K70").Value = 8367.846
$ws.Range("L70").Value = 9943
$ws.Range("M70").Value = -8097.846
$ws.Range("N70").Value = -10483
$ws.Range("H73").Value = 8919.15
$ws.Range("I73").Value = 8367.846
$ws.Range("J73").Value = 9943
$ws.Range("K73").Value = 8367.846
$ws.Range("L73").Value = 9943
$ws.Range("M73").Value = -7431.846
$ws.Range("N73").Value = -11815

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 20003
$ws.Range("I10").Value = 20003
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 20003
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -19863
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H16").Value = 199.52942
$ws.Range("I16").Value = 205.14285
$ws.Range("J16").Value = 173.33333
$ws.Range("K16").Value = 205.14285
$ws.Range("L16").Value = 173.33333
$ws.Range("M16").Value = -35.14285000000001
$ws.Range("N16").Value = -513.3333299999999
$ws.Range("H46").Value = 2309.2104
$ws.Range("I46").Value = 1388.8
$ws.Range("J46").Value = 3331.889
$ws.Range("K46").Value = 1388.8
$ws.Range("L46").Value = 3331.889
$ws.Range("M46").Value = -1200.8
$ws.Range("N46").Value = -3707.889
$ws.Range("H93").Value = 4118.5
$ws.Range("I93").Value = 4145.9
$ws.Range("J93").Value = 4050
$ws.Range("K93").Value = 4145.9
$ws.Range("L93").Value = 4050
$ws.Range("M93").Value = -2897.9
$ws.Range("N93").Value = -6546
$ws.Range("H122").Value = 4261.788
$ws.Range("I122").Value = 3952.9092
$ws.Range("J122").Value = 4879.5454
$ws.Range("K122").Value = 11858.7276
$ws.Range("L122").Value = 14638.6362
$ws.Range("M122").Value = -9408.7276
$ws.Range("N122").Value = -19538.6362
$ws.Range("H132").Value = 791576.5
$ws.Range("I132").Value = 1072639.8
$ws.Range("J132").Value = 4599.4
$ws.Range("K132").Value = 3217919.4
$ws.Range("L132").Value = 13798.2
$ws.Range("M132").Value = -3215389.4
$ws.Range("N132").Value = -18858.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 1048.95
$ws.Range("I4").Value = 780.125
$ws.Range("J4").Value = 2124.25
$ws.Range("K4").Value = 780.125
$ws.Range("L4").Value = 2124.25
$ws.Range("M4").Value = -667.125
$ws.Range("N4").Value = -2350.25
$ws.Range("H132").Value = 15944.814
$ws.Range("I132").Value = 16432.863
$ws.Range("J132").Value = 13797.4
$ws.Range("K132").Value = 49298.58900000001
$ws.Range("L132").Value = 41392.2
$ws.Range("M132").Value = -46768.58900000001
$ws.Range("N132").Value = -46452.2
